# "Generate Report for Handoff"
#
# The f8a331d6-... localization entry moved from "In Translation" to
# "Ready for handoff": status text + handoff datetime are refreshed on the
# Overview sheet and on each per-locale (zh-cn / de-de) sheet, and the
# zh-cn/de-de "Priority" for that row flips from "ht" to "mt" with a new
# handoff timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the f8a331d6-... file
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-14 00:19:38"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is the f8a331d6-... file
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-14 00:19:30"

# ---------------------------------------------------------------------
# de-de sheet - row 3 is the f8a331d6-... file
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-14 00:19:38"

# ---------------------------------------------------------------------
# The longer "Ready for handoff" status text widened the Status columns
# (Overview!E:F and the "Status" column on each locale sheet).
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
